# Applies "some adjustments to lesson1" to the language workbook.
#
# Net effect on the Key/Value table (sheet "en"):
#   - lesson1_distribute_1 value text is reworded
#   - two brand new rows are inserted right after it:
#       lesson1_distribute_2, lesson1_distribute_3
#   - lesson1_area_1 / lesson1_area_2 / lesson1_area_3 keep their keys but
#     get new value text
#   - lesson1_area_4 is removed entirely
#
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows right after row 44 (lesson1_distribute_1) to make
# room for the new lesson1_distribute_2 / lesson1_distribute_3 entries.
$ws.Rows("45:46").Insert()

# Reword the existing lesson1_distribute_1 value.
$ws.Range("B44").Value = "By using the distributive property, we can split a large number up into smaller pieces."

# New lesson1_distribute_2 row.
$ws.Range("A45").Value = "lesson1_distribute_2"
$ws.Range("B45").Value = "In this example, we split up 12 into: 10 and 2. Making it easier to solve the equation."

# New lesson1_distribute_3 row.
$ws.Range("A46").Value = "lesson1_distribute_3"
$ws.Range("B46").Value = "Pay close attention to how the distribution works!"

# The old lesson1_area_1/_2/_3 rows (now shifted down to 47/48/49) get new
# wording while keeping their original keys.
$ws.Range("B47").Value = "To help us visualize this further, consider the product as the area of a rectangle."
$ws.Range("B48").Value = "Splitting the rectangle into two is much similar to how the distributive property works."
$ws.Range("B49").Value = "You compute the area of the two rectangles, and then add them up to get the total area."

# lesson1_area_4 (now shifted down to row 50) is dropped entirely.
$ws.Rows("50:50").Delete()

# Reflect the author's final scroll position/selection in the saved file.
$ws.Application.ActiveWindow.ScrollRow = 34
$ws.Range("A49").Select()
